$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.925.43"
Set-TextValue $ws.Range("E2") "  -6.58%  "
Set-TextValue $ws.Range("D3") "2.431.35"
Set-TextValue $ws.Range("E3") "  -9.40%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "528.36"
Set-TextValue $ws.Range("E5") "  -3.82%  "
Set-TextValue $ws.Range("D6") "146.17"
Set-TextValue $ws.Range("E6") "  -7.20%  "
Set-TextValue $ws.Range("D7") "0.998"
Set-TextValue $ws.Range("E7") "  -0.10%  "
Set-TextValue $ws.Range("D8") "0.564"
Set-TextValue $ws.Range("E8") "  -3.97%  "
Set-TextValue $ws.Range("E9") "  -6.81%  "
Set-TextValue $ws.Range("E10") "  -2.35%  "
Set-TextValue $ws.Range("D11") "5.32"
Set-TextValue $ws.Range("E11") "  +4.18%  "
Set-TextValue $ws.Range("E12") "  -5.40%  "
Set-TextValue $ws.Range("D13") "2.867.05"
Set-TextValue $ws.Range("E13") "  -9.18%  "
Set-TextValue $ws.Range("D14") "23.99"
Set-TextValue $ws.Range("E14") "  -7.93%  "
Set-TextValue $ws.Range("D15") "58.858.97"
Set-TextValue $ws.Range("E15") "  -6.49%  "
Set-TextValue $ws.Range("E16") "  -6.80%  "
Set-TextValue $ws.Range("D17") "2.473.45"
Set-TextValue $ws.Range("E17") "  -7.84%  "
Set-TextValue $ws.Range("D18") "11.03"
Set-TextValue $ws.Range("E18") "  -7.86%  "
Set-TextValue $ws.Range("D19") "4.29"
Set-TextValue $ws.Range("E19") "  -6.04%  "
Set-TextValue $ws.Range("D20") "322.14"
Set-TextValue $ws.Range("E20") "  -6.21%  "
Set-TextValue $ws.Range("D21") "0.969"
Set-TextValue $ws.Range("E21") "  -2.96%  "
Set-TextValue $ws.Range("D22") "5.73"
Set-TextValue $ws.Range("E22") "  -9.19%  "
Set-TextValue $ws.Range("D23") "0.465"
Set-TextValue $ws.Range("E23") "  -7.62%  "
Set-TextValue $ws.Range("D24") "60.29"
Set-TextValue $ws.Range("E24") "  -4.90%  "
Set-TextValue $ws.Range("D25") "0.160"
Set-TextValue $ws.Range("E25") "  -4.49%  "
Set-TextValue $ws.Range("D26") "0.978"
Set-TextValue $ws.Range("E26") "  -2.26%  "
Set-TextValue $ws.Range("D27") "7.73"
Set-TextValue $ws.Range("E27") "  -5.03%  "
Set-TextValue $ws.Range("D28") "6.85"
Set-TextValue $ws.Range("E28") "  -2.22%  "
Set-TextValue $ws.Range("D29") "1.27"
Set-TextValue $ws.Range("E29") "  -5.04%  "
Set-TextValue $ws.Range("E30") "  -6.18%  "
Set-TextValue $ws.Range("D31") "0.0₃0770"
Set-TextValue $ws.Range("E31") "  -9.55%  "
Set-TextValue $ws.Range("D32") "0.998"
Set-TextValue $ws.Range("E32") "  -0.07%  "
Set-TextValue $ws.Range("D33") "157.72"
Set-TextValue $ws.Range("E33") "  -4.46%  "
Set-TextValue $ws.Range("D34") "4.49"
Set-TextValue $ws.Range("E34") "  -6.77%  "
Set-TextValue $ws.Range("D35") "18.27"
Set-TextValue $ws.Range("E35") "  -6.36%  "
Set-TextValue $ws.Range("E36") "  -6.88%  "
Set-TextValue $ws.Range("D37") "1.71"
Set-TextValue $ws.Range("E37") "  -3.75%  "
Set-TextValue $ws.Range("D38") "5.71"
Set-TextValue $ws.Range("E38") "  -7.30%  "
Set-TextValue $ws.Range("D39") "308.12"
Set-TextValue $ws.Range("E39") "  -9.28%  "
Set-TextValue $ws.Range("B40") "OKB"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D40") "36.62"
Set-TextValue $ws.Range("E40") "  -4.03%  "
Set-TextValue $ws.Range("B41") "SuiNetwork"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D41") "0.842"
Set-TextValue $ws.Range("E41") "  -9.66%  "
Set-TextValue $ws.Range("D42") "3.72"
Set-TextValue $ws.Range("E42") "  -5.42%  "
Set-TextValue $ws.Range("D43") "0.998"
Set-TextValue $ws.Range("E43") "  -0.07%  "
Set-TextValue $ws.Range("D44") "10.70"
Set-TextValue $ws.Range("E44") "  -3.08%  "
Set-TextValue $ws.Range("B45") "Mantle"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D45") "0.576"
Set-TextValue $ws.Range("E45") "  -6.79%  "
Set-TextValue $ws.Range("B46") "Stellar"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D46") "0.0930"
Set-TextValue $ws.Range("E46") "  -4.35%  "
Set-TextValue $ws.Range("D47") "19.11"
Set-TextValue $ws.Range("E47") "  -7.81%  "
Set-TextValue $ws.Range("E48") "  -7.38%  "
Set-TextValue $ws.Range("D49") "18.47"
Set-TextValue $ws.Range("E49") "  -9.19%  "
Set-TextValue $ws.Range("D50") "1.981.48"
Set-TextValue $ws.Range("E50") "  -5.07%  "
Set-TextValue $ws.Range("D51") "0.0227"
Set-TextValue $ws.Range("E51") "  -5.31%  "
